$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G3").Value = 1.83
$ws.Range("I3").Value = 4.6
$ws.Range("J3").Value = 2.4
$ws.Range("K3").Value = 2.05
$ws.Range("L3").Value = 4.9
$ws.Range("P3").Value = 2.7
$ws.Range("Q3").Value = 2.22
$ws.Range("T3").Value = 2.57
$ws.Range("W3").Value = 5.9
$ws.Range("X3").Value = 8
$ws.Range("Z3").Value = 15.5
$ws.Range("AA3").Value = 16
$ws.Range("AD3").Value = 6.1
$ws.Range("AE3").Value = 16.5
$ws.Range("AH3").Value = 10.25
$ws.Range("AI3").Value = 25
$ws.Range("AJ3").Value = 15.5
$ws.Range("AK3").Value = 90
$ws.Range("AL3").Value = 55
$ws.Range("AM3").Value = 60
$ws.Range("AN3").Value = 3.6
$ws.Range("AO3").Value = 9.25
$ws.Range("AT3").Value = 2.57
$ws.Range("AV3").Value = 70
$ws.Range("AW3").Value = 6.3
$ws.Range("AX3").Value = 28
$ws.Range("AY3").Value = 35
$ws.Range("AZ3").Value = 175

$ws.Range("L5").Value = 7
$ws.Range("M5").Value = 1.03
$ws.Range("N5").Value = 10
$ws.Range("U5").Value = 2.1
$ws.Range("V5").Value = 1.67
$ws.Range("W5").Value = 7
$ws.Range("X5").Value = 6.5
$ws.Range("AF5").Value = 67
$ws.Range("AJ5").Value = 21
$ws.Range("AU5").Value = 9.5

$ws.Range("M6").Value = 1.06
$ws.Range("N6").Value = 8
